# Generate Report for Handback
#
# Fills in the "Latest Target File" (I), "Latest Handback File" (J),
# "Latest Handback DateTime" (K) and "Error Detail" (P) columns for row 7
# (the d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.md file) on both the "zh-cn"
# and "de-de" worksheets, now that a handback report has been generated
# for that file.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd51bae0fe2a468fbdfcc0f3cf174999cb3dd00a/e2e/d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.md"
$targetDisplay = "d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8fcb1c575c7badc01e8697af8dcd667743573ef/e2e/d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd51bae0fe2a468fbdfcc0f3cf174999cb3dd00a/e2e/d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.md."

# --- zh-cn sheet, row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.8cf7817f5f743a9a89ebac3b609edf354e9ce14f.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-17 02:57:37"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetDisplay)

# --- de-de sheet, row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "d3c3f1a0-43c3-48f7-afb3-5a495624c8d6.8cf7817f5f743a9a89ebac3b609edf354e9ce14f.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-17 02:57:44"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetDisplay)
